$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-07 05:18:01'
$ws.Range('N2').Value = '-2.7 °C 4:49 TU'
$ws.Range('O2').Value = '-1.4 °C'
$ws.Range('E3').Value = '2026-02-07 05:18:04'
$ws.Range('H3').Value = '''94%'
$ws.Range('N3').Value = '-7.0 °C 4:43 TU'
$ws.Range('O3').Value = '-5.7 °C'
$ws.Range('E4').Value = '2026-02-07 05:18:07'
$ws.Range('H4').Value = '''57%'
$ws.Range('J4').Value = '1001.0 hPa'
$ws.Range('O4').Value = '11.5 °C'
$ws.Range('E5').Value = '2026-02-07 05:18:09'
$ws.Range('N5').Value = '7.6 °C 4:59 TU'
$ws.Range('O5').Value = '8.9 °C'
$ws.Range('E6').Value = '2026-02-07 05:18:12'
$ws.Range('J6').Value = '1002.8 hPa'
$ws.Range('N6').Value = '11.3 °C 4:38 TU'
$ws.Range('E7').Value = '2026-02-07 05:18:15'
$ws.Range('H7').Value = '''76%'
$ws.Range('J7').Value = '1002.5 hPa'
$ws.Range('N7').Value = '7.0 °C 4:44 TU'
$ws.Range('O7').Value = '7.7 °C'
$ws.Range('E8').Value = '2026-02-07 05:18:17'
$ws.Range('N8').Value = '2.6 °C 4:53 TU'
$ws.Range('O8').Value = '4.2 °C'
$ws.Range('E9').Value = '2026-02-07 05:18:20'
$ws.Range('N9').Value = '-0.2 °C 4:59 TU'
$ws.Range('O9').Value = '1.9 °C'
$ws.Range('E10').Value = '2026-02-07 05:18:23'
$ws.Range('E11').Value = '2026-02-07 05:18:25'
$ws.Range('J11').Value = '1005.4 hPa'
$ws.Range('E12').Value = '2026-02-07 05:18:28'
$ws.Range('H12').Value = '''72%'
$ws.Range('N12').Value = '7.0 °C 4:59 TU'
$ws.Range('O12').Value = '9.7 °C'
$ws.Range('E13').Value = '2026-02-07 05:18:31'
$ws.Range('M13').Value = '9.9 °C 4:59 TU'
$ws.Range('O13').Value = '7.2 °C'
$ws.Range('E14').Value = '2026-02-07 05:18:34'
$ws.Range('H14').Value = '''77%'
$ws.Range('E15').Value = '2026-02-07 05:18:36'
$ws.Range('H15').Value = '''84%'
$ws.Range('J15').Value = '1001.4 hPa'
$ws.Range('N15').Value = '3.8 °C 4:45 TU'
$ws.Range('O15').Value = '6.6 °C'
$ws.Range('E16').Value = '2026-02-07 05:18:39'
$ws.Range('L16').Value = '18.7 km/h - 294º 4:52 TU'
$ws.Range('N16').Value = '1.8 °C 4:34 TU'
$ws.Range('O16').Value = '3.0 °C'
$ws.Range('E17').Value = '2026-02-07 05:18:42'
$ws.Range('J17').Value = '1004.6 hPa'
$ws.Range('N17').Value = '2.8 °C 4:54 TU'
$ws.Range('E18').Value = '2026-02-07 05:18:44'
$ws.Range('N18').Value = '-8.9 °C 4:59 TU'
$ws.Range('O18').Value = '-7.2 °C'
$ws.Range('E19').Value = '2026-02-07 05:18:47'
$ws.Range('I19').Value = '0.2 mm'
$ws.Range('J19').Value = '1005.9 hPa'
$ws.Range('O19').Value = '4.6 °C'
$ws.Range('E20').Value = '2026-02-07 05:18:50'
$ws.Range('H20').Value = '''84%'
$ws.Range('E21').Value = '2026-02-07 05:18:52'
$ws.Range('H21').Value = '''74%'
$ws.Range('J21').Value = '1001.6 hPa'
$ws.Range('N21').Value = '2.8 °C 4:52 TU'
$ws.Range('O21').Value = '7.0 °C'
$ws.Range('E22').Value = '2026-02-07 05:18:55'
$ws.Range('H22').Value = '''93%'
$ws.Range('M22').Value = '8.0 °C 4:55 TU'
$ws.Range('O22').Value = '5.9 °C'
$ws.Range('E23').Value = '2026-02-07 05:18:58'
$ws.Range('J23').Value = '1001.3 hPa'
$ws.Range('L23').Value = '20.9 km/h - 312º 4:45 TU'
$ws.Range('E24').Value = '2026-02-07 05:19:00'
$ws.Range('J24').Value = '1000.6 hPa'
$ws.Range('E25').Value = '2026-02-07 05:19:03'
$ws.Range('J25').Value = '1005.1 hPa'
$ws.Range('E26').Value = '2026-02-07 05:19:06'
$ws.Range('N26').Value = '-4.3 °C 4:53 TU'
$ws.Range('O26').Value = '-1.7 °C'
$ws.Range('E27').Value = '2026-02-07 05:19:08'
$ws.Range('J27').Value = '1001.1 hPa'
$ws.Range('E28').Value = '2026-02-07 05:19:11'
$ws.Range('J28').Value = '1003.8 hPa'
$ws.Range('N28').Value = '1.7 °C 4:56 TU'
$ws.Range('O28').Value = '3.3 °C'
$ws.Range('E29').Value = '2026-02-07 05:19:13'
$ws.Range('H29').Value = '''58%'
$ws.Range('O29').Value = '11.1 °C'
$ws.Range('E30').Value = '2026-02-07 05:19:16'
$ws.Range('H30').Value = '''83%'
$ws.Range('E31').Value = '2026-02-07 05:19:19'
$ws.Range('J31').Value = '1005.7 hPa'
$ws.Range('N31').Value = '3.4 °C 4:51 TU'
$ws.Range('E32').Value = '2026-02-07 05:19:21'
$ws.Range('H32').Value = '''60%'
$ws.Range('J32').Value = '1004.2 hPa'
$ws.Range('L32').Value = '33.1 km/h - 283º 4:55 TU'
$ws.Range('E33').Value = '2026-02-07 05:19:24'
$ws.Range('H33').Value = '''90%'
$ws.Range('N33').Value = '5.3 °C 4:40 TU'
$ws.Range('O33').Value = '7.1 °C'
$ws.Range('E34').Value = '2026-02-07 05:19:26'
$ws.Range('O34').Value = '6.5 °C'
$ws.Range('E35').Value = '2026-02-07 05:19:29'
$ws.Range('N35').Value = '-8.3 °C 4:32 TU'
$ws.Range('O35').Value = '-5.4 °C'
$ws.Range('E36').Value = '2026-02-07 05:19:32'
$ws.Range('I36').Value = '0.1 mm'
$ws.Range('J36').Value = '1006.4 hPa'
$ws.Range('O36').Value = '4.7 °C'
